$d = $word.ActiveDocument

$d.Content.Find.Execute("2024-03-02 Saturday", $true, $false, $false, $false, $false, $true, 1, $false, "2024-03-03 Sunday", 2) | Out-Null
$d.Content.Find.Execute("20÷3=6, 2", $true, $false, $false, $false, $false, $true, 1, $false, "14÷5=2, 4", 2) | Out-Null
$d.Content.Find.Execute("25÷2=12, 1", $true, $false, $false, $false, $false, $true, 1, $false, "93÷5=18, 3", 2) | Out-Null
$d.Content.Find.Execute("95÷9=10, 5", $true, $false, $false, $false, $false, $true, 1, $false, "42÷3=14, 0", 2) | Out-Null
$d.Content.Find.Execute("93÷7=13, 2", $true, $false, $false, $false, $false, $true, 1, $false, "96÷8=12, 0", 2) | Out-Null
$d.Content.Find.Execute("80÷6=13, 2", $true, $false, $false, $false, $false, $true, 1, $false, "12÷7=1, 5", 2) | Out-Null
$d.Content.Find.Execute("12÷8=1, 4", $true, $false, $false, $false, $false, $true, 1, $false, "61÷6=10, 1", 2) | Out-Null
$d.Content.Find.Execute("39÷9=4, 3", $true, $false, $false, $false, $false, $true, 1, $false, "56÷4=14, 0", 2) | Out-Null
$d.Content.Find.Execute("22÷9=2, 4", $true, $false, $false, $false, $false, $true, 1, $false, "60÷8=7, 4", 2) | Out-Null
$d.Content.Find.Execute("61÷4=15, 1", $true, $false, $false, $false, $false, $true, 1, $false, "20÷5=4, 0", 2) | Out-Null
$d.Content.Find.Execute("16÷4=4, 0", $true, $false, $false, $false, $false, $true, 1, $false, "42÷7=6, 0", 2) | Out-Null
$d.Content.Find.Execute("37÷4=9, 1", $true, $false, $false, $false, $false, $true, 1, $false, "57÷8=7, 1", 2) | Out-Null
$d.Content.Find.Execute("97÷3=32, 1", $true, $false, $false, $false, $false, $true, 1, $false, "70÷4=17, 2", 2) | Out-Null
$d.Content.Find.Execute("43÷3=14, 1", $true, $false, $false, $false, $false, $true, 1, $false, "56÷6=9, 2", 2) | Out-Null
$d.Content.Find.Execute("64÷2=32, 0", $true, $false, $false, $false, $false, $true, 1, $false, "96÷5=19, 1", 2) | Out-Null
$d.Content.Find.Execute("90÷4=22, 2", $true, $false, $false, $false, $false, $true, 1, $false, "67÷3=22, 1", 2) | Out-Null
$d.Content.Find.Execute("78÷7=11, 1", $true, $false, $false, $false, $false, $true, 1, $false, "58÷7=8, 2", 2) | Out-Null
$d.Content.Find.Execute("22÷5=4, 2", $true, $false, $false, $false, $false, $true, 1, $false, "50÷7=7, 1", 2) | Out-Null
$d.Content.Find.Execute("93÷3=31, 0", $true, $false, $false, $false, $false, $true, 1, $false, "29÷5=5, 4", 2) | Out-Null
$d.Content.Find.Execute("52÷3=17, 1", $true, $false, $false, $false, $false, $true, 1, $false, "47÷9=5, 2", 2) | Out-Null
$d.Content.Find.Execute("48÷5=9, 3", $true, $false, $false, $false, $false, $true, 1, $false, "46÷6=7, 4", 2) | Out-Null
$d.Content.Find.Execute("97÷9=10, 7", $true, $false, $false, $false, $false, $true, 1, $false, "44÷5=8, 4", 2) | Out-Null
$d.Content.Find.Execute("83÷8=10, 3", $true, $false, $false, $false, $false, $true, 1, $false, "60÷4=15, 0", 2) | Out-Null
$d.Content.Find.Execute("79÷7=11, 2", $true, $false, $false, $false, $false, $true, 1, $false, "34÷3=11, 1", 2) | Out-Null
$d.Content.Find.Execute("12÷5=2, 2", $true, $false, $false, $false, $false, $true, 1, $false, "49÷4=12, 1", 2) | Out-Null
$d.Content.Find.Execute("87÷8=10, 7", $true, $false, $false, $false, $false, $true, 1, $false, "75÷6=12, 3", 2) | Out-Null
